# Apply updated crypto price/volume data to the worksheet.
# Numeric-looking Price values in column D must stay as TEXT (as in the
# source data), so they are written with a leading apostrophe just like
# typing them directly into Excel would force a text entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''56.575.32'
$ws.Range("E2").Value = '  +2.62%  '
$ws.Range("D3").Value = '''2.320.98'
$ws.Range("E3").Value = '  +2.23%  '
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").Value = '''517.19'
$ws.Range("E5").Value = '  +2.46%  '
$ws.Range("D6").Value = '''135.32'
$ws.Range("E6").Value = '  +5.80%  '
$ws.Range("E7").Value = '  -0.15%  '
$ws.Range("D8").Value = '''0.537'
$ws.Range("E8").Value = '  +1.52%  '
$ws.Range("D9").Value = '''2.339.41'
$ws.Range("E9").Value = '  +2.67%  '
$ws.Range("E10").Value = '  +4.19%  '
$ws.Range("E11").Value = '  -1.09%  '
$ws.Range("E12").Value = '  +5.17%  '
$ws.Range("D13").Value = '''0.342'
$ws.Range("E13").Value = '  +0.55%  '
$ws.Range("D14").Value = '''23.98'
$ws.Range("E14").Value = '  +2.01%  '
$ws.Range("D15").Value = '''2.732.78'
$ws.Range("E15").Value = '  +2.23%  '
$ws.Range("D16").Value = '''56.638.48'
$ws.Range("E16").Value = '  +3.00%  '
$ws.Range("E17").Value = '  +2.90%  '
$ws.Range("D18").Value = '''2.338.88'
$ws.Range("E18").Value = '  +2.78%  '
$ws.Range("D19").Value = '''10.53'
$ws.Range("E19").Value = '  +1.81%  '
$ws.Range("E20").Value = '  +0.95%  '
$ws.Range("D21").Value = '''324.69'
$ws.Range("E21").Value = '  +3.70%  '
$ws.Range("E22").Value = '  +0.30%  '
$ws.Range("E23").Value = '  +0.18%  '
$ws.Range("D24").Value = '''60.72'
$ws.Range("E24").Value = '  +1.40%  '
$ws.Range("D25").Value = '''0.165'
$ws.Range("E25").Value = '  +6.88%  '
$ws.Range("D26").Value = '''0.993'
$ws.Range("E26").Value = '  -0.30%  '
$ws.Range("D27").Value = '''7.97'
$ws.Range("E27").Value = '  +6.09%  '
$ws.Range("D28").Value = '''1.28'
$ws.Range("E28").Value = '  +11.96%  '
$ws.Range("D29").Value = '''0.0₃0739'
$ws.Range("E29").Value = '  +5.18%  '
$ws.Range("D30").Value = '''167.43'
$ws.Range("E30").Value = '  -2.05%  '
$ws.Range("D31").Value = '''1.69'
$ws.Range("E31").Value = '  +3.44%  '
$ws.Range("D32").Value = '''6.19'
$ws.Range("E32").Value = '  +0.72%  '
$ws.Range("D33").Value = '''18.46'
$ws.Range("E33").Value = '  +2.83%  '
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("D35").Value = '''0.992'
$ws.Range("E35").Value = '  -0.43%  '
$ws.Range("D36").Value = '''1.26'
$ws.Range("E36").Value = '  +2.59%  '
$ws.Range("D37").Value = '''0.919'
$ws.Range("E37").Value = '  +1.89%  '
$ws.Range("D38").Value = '''4.01'
$ws.Range("E38").Value = '  +3.23%  '
$ws.Range("E39").Value = '  +7.15%  '
$ws.Range("D40").Value = '''38.29'
$ws.Range("E40").Value = '  +4.51%  '
$ws.Range("D41").Value = '''0.380'
$ws.Range("E41").Value = '  +1.86%  '
$ws.Range("D42").Value = '''141.42'
$ws.Range("E42").Value = '  +4.28%  '
$ws.Range("E43").Value = '  +3.75%  '
$ws.Range("E44").Value = '  +7.54%  '
$ws.Range("D45").Value = '''276.07'
$ws.Range("E45").Value = '  +7.19%  '
$ws.Range("E46").Value = '  +1.97%  '
$ws.Range("D47").Value = '''0.0507'
$ws.Range("E47").Value = '  +0.61%  '
$ws.Range("D48").Value = '''0.561'
$ws.Range("E48").Value = '  +2.78%  '
$ws.Range("D49").Value = '''0.0219'
$ws.Range("E49").Value = '  +3.12%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").Value = '''17.83'
$ws.Range("E50").Value = '  +9.07%  '
$ws.Range("B51").Value = 'Polygon'
$ws.Range("C51").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D51").Value = '''0.380'
$ws.Range("E51").Value = '  +1.69%  '
